$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 4474.75
$ws.Range("I40").Value = 3437.25
$ws.Range("J40").Value = 6549.75
$ws.Range("K40").Value = 3437.25
$ws.Range("L40").Value = 6549.75
$ws.Range("M40").Value = -3262.25
$ws.Range("N40").Value = -6899.75

# Row 62
$ws.Range("H62").Value = 83385500
$ws.Range("I62").Value = 166682130
$ws.Range("J62").Value = 88878
$ws.Range("K62").Value = 166682130
$ws.Range("L62").Value = 88878
$ws.Range("M62").Value = -166681506
$ws.Range("N62").Value = -90126

# Row 64
$ws.Range("H64").Value = 5547.75
$ws.Range("I64").Value = 4002
$ws.Range("J64").Value = 6063
$ws.Range("K64").Value = 4002
$ws.Range("L64").Value = 6063
$ws.Range("M64").Value = -3754
$ws.Range("N64").Value = -6559

# Row 65
$ws.Range("H65").Value = 83385500
$ws.Range("I65").Value = 166682130
$ws.Range("J65").Value = 88878
$ws.Range("K65").Value = 833410650
$ws.Range("L65").Value = 444390
$ws.Range("M65").Value = -833407530
$ws.Range("N65").Value = -450630

# Row 67
$ws.Range("H67").Value = 5547.75
$ws.Range("I67").Value = 4002
$ws.Range("J67").Value = 6063
$ws.Range("K67").Value = 4002
$ws.Range("L67").Value = 6063
$ws.Range("M67").Value = -3144
$ws.Range("N67").Value = -7779

# Row 106
$ws.Range("H106").Value = 1056.1428
$ws.Range("I106").Value = 1056.1428
$ws.Range("K106").Value = 1056.1428
$ws.Range("M106").Value = -425.1428000000001

# Row 123
$ws.Range("H123").Value = 51598
$ws.Range("J123").Value = 51598
$ws.Range("L123").Value = 51598
$ws.Range("N123").Value = -61398

# Row 132
$ws.Range("H132").Value = 1213.6123
$ws.Range("I132").Value = 1222.7234
$ws.Range("K132").Value = 3668.1702
$ws.Range("M132").Value = -1138.1702

# Row 137
$ws.Range("H137").Value = 2700.4856
$ws.Range("I137").Value = 2418.1738
$ws.Range("J137").Value = 3241.5833
$ws.Range("K137").Value = 7254.5214
$ws.Range("L137").Value = 9724.749899999999
$ws.Range("M137").Value = -4704.5214
$ws.Range("N137").Value = -14824.7499

# Row 141
$ws.Range("H141").Value = 13336065
$ws.Range("I141").Value = 15154038
$ws.Range("K141").Value = 45462114
$ws.Range("M141").Value = -45456934

$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Range("H22").Value = 812.5
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -201
$ws.Range("N22").Value = -1598

# Row 32
$ws.Range("H32").Value = 1898374.2
$ws.Range("I32").Value = 2053610.8
$ws.Range("K32").Value = 2053610.8
$ws.Range("M32").Value = -2053323.8

# Row 45
$ws.Range("H45").Value = 7627.222
$ws.Range("I45").Value = 3050.5
$ws.Range("K45").Value = 3050.5
$ws.Range("M45").Value = -2673.5

# Row 74
$ws.Range("H74").Value = 34065.215
$ws.Range("I74").Value = 54745.6
$ws.Range("K74").Value = 54745.6
$ws.Range("M74").Value = -53871.6

# Row 77
$ws.Range("H77").Value = 34065.215
$ws.Range("I77").Value = 54745.6
$ws.Range("K77").Value = 273728
$ws.Range("M77").Value = -269360

# Row 102
$ws.Range("H102").Value = 989
$ws.Range("I102").Value = 988.8946999999999
$ws.Range("K102").Value = 988.8946999999999
$ws.Range("M102").Value = 633.1053000000001

# Row 132
$ws.Range("H132").Value = 5395.396
$ws.Range("I132").Value = 4349.5806
$ws.Range("J132").Value = 7302.4707
$ws.Range("K132").Value = 13048.7418
$ws.Range("L132").Value = 21907.4121
$ws.Range("M132").Value = -10518.7418
$ws.Range("N132").Value = -26967.4121

$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 5000
$ws.Range("I5").Value = 5000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 5000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -4887
$ws.Range("N5").ClearContents()

# Row 19
$ws.Range("H19").Value = 7000
$ws.Range("I19").Value = 7000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 7000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -6827
$ws.Range("N19").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7319.357
$ws.Range("I31").Value = 2810.2354
$ws.Range("J31").Value = 10385.56
$ws.Range("K31").Value = 2810.2354
$ws.Range("L31").Value = 10385.56
$ws.Range("M31").Value = -2515.2354
$ws.Range("N31").Value = -10975.56

# Row 34
$ws.Range("H34").Value = 7319.357
$ws.Range("I34").Value = 2810.2354
$ws.Range("J34").Value = 10385.56
$ws.Range("K34").Value = 2810.2354
$ws.Range("L34").Value = 10385.56
$ws.Range("M34").Value = -2608.2354
$ws.Range("N34").Value = -10789.56

# Row 64
$ws.Range("H64").Value = 40600
$ws.Range("J64").Value = 40600
$ws.Range("L64").Value = 40600
$ws.Range("N64").Value = -41096

# Row 67
$ws.Range("H67").Value = 40600
$ws.Range("J67").Value = 40600
$ws.Range("L67").Value = 40600
$ws.Range("N67").Value = -42316

# Row 99
$ws.Range("H99").Value = 8126.778
$ws.Range("I99").Value = 9471.875
$ws.Range("K99").Value = 9471.875
$ws.Range("M99").Value = -7973.875

# Row 126
$ws.Range("H126").Value = 8126.778
$ws.Range("I126").Value = 9471.875
$ws.Range("K126").Value = 28415.625
$ws.Range("M126").Value = -25945.625

# Row 132
$ws.Range("H132").Value = 12909878
$ws.Range("I132").Value = 2113.9167
$ws.Range("J132").Value = 21062150
$ws.Range("K132").Value = 6341.750100000001
$ws.Range("L132").Value = 63186450
$ws.Range("M132").Value = -3811.750100000001
$ws.Range("N132").Value = -63191510

$ws = $wb.Worksheets.Item("CUL")
# Row 36
$ws.Range("H36").Value = 202010160
$ws.Range("I36").Value = 250012700
$ws.Range("K36").Value = 750038100
$ws.Range("M36").Value = -750037931

# Row 121
$ws.Range("H121").Value = 913
$ws.Range("J121").Value = 1342.8334
$ws.Range("L121").Value = 4028.5002
$ws.Range("N121").Value = -6648.5002

# Row 122
$ws.Range("H122").Value = 2831447
$ws.Range("J122").Value = 4038.6
$ws.Range("L122").Value = 36347.4
$ws.Range("N122").Value = -41247.4

# Row 134
$ws.Range("H134").Value = 51905.715
$ws.Range("I134").Value = 59333.777
$ws.Range("J134").Value = 7337.3335
$ws.Range("K134").Value = 178001.331
$ws.Range("L134").Value = 22012.0005
$ws.Range("M134").Value = -172931.331
$ws.Range("N134").Value = -32152.0005

# Row 138
$ws.Range("H138").Value = 94355.27
$ws.Range("I138").Value = 102887.9
$ws.Range("K138").Value = 308663.7
$ws.Range("M138").Value = -303523.7

# Row 139
$ws.Range("H139").Value = 44036.152
$ws.Range("I139").Value = 69032.734
$ws.Range("J139").Value = 9949.909
$ws.Range("K139").Value = 207098.202
$ws.Range("L139").Value = 29849.727
$ws.Range("M139").Value = -201958.202
$ws.Range("N139").Value = -40129.727

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 667224.9399999999
$ws.Range("I107").Value = 1000275
$ws.Range("K107").Value = 1000275
$ws.Range("M107").Value = -998355

# Row 132
$ws.Range("H132").Value = 4154.393
$ws.Range("I132").Value = 1444.1111
$ws.Range("K132").Value = 4332.3333
$ws.Range("M132").Value = -1802.3333

$ws = $wb.Worksheets.Item("LTW")
# Row 36
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

# Row 132
$ws.Range("H132").Value = 10644919
$ws.Range("I132").Value = 25002890
$ws.Range("J132").Value = 9385.519
$ws.Range("K132").Value = 75008670
$ws.Range("L132").Value = 28156.557
$ws.Range("M132").Value = -75006140
$ws.Range("N132").Value = -33216.557

# Row 136
$ws.Range("H136").Value = 8628.1875
$ws.Range("I136").Value = 1137.5454
$ws.Range("J136").Value = 12551.857
$ws.Range("K136").Value = 3412.6362
$ws.Range("L136").Value = 37655.571
$ws.Range("M136").Value = -862.6361999999999
$ws.Range("N136").Value = -42755.571

$ws = $wb.Worksheets.Item("WVR")
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Row 132
$ws.Range("H132").Value = 14717814
$ws.Range("I132").Value = 18522880
$ws.Range("K132").Value = 55568640
$ws.Range("M132").Value = -55566110
